# The Stage Frights - apply edit described by the diff:
#  1. "merch.html THE SWAG PAGE" -> "swag.html THE SWAG PAGE"
#     (i.e. replace "merch" with "swag" in that paragraph)
#  2. Word auto-tracks the last edit location with a hidden "_GoBack"
#     bookmark. Previously it sat at the end of the "Tour Name" paragraph;
#     after this edit it should instead sit right after the inserted
#     "swag" text (and get removed from its old spot, since a document
#     only ever has a single "_GoBack" bookmark).

$d = $word.ActiveDocument

# Step 1: replace "merch" with "swag" (only the occurrence in
# "merch.html THE SWAG PAGE")
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("merch.html THE SWAG PAGE", $true, $false, $false, $false, $false, $true, 1, $false, "swag.html THE SWAG PAGE", 2) | Out-Null

# Step 2: move the "_GoBack" bookmark to sit right after the newly
# inserted "swag" text. Adding a bookmark named "_GoBack" replaces any
# existing bookmark of that name (Word keeps only one), which also takes
# care of removing it from the end of the "Tour Name" paragraph.
$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Execute("swag") | Out-Null
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r) | Out-Null
